# Apply the edits described by the diff:
# - Update Case 4 - study 1 (pc) block:
#     [nx, nz]  : [50, 100]  -> [75, 75]
#     [lx, lz]  : [600, 600] -> [800, 450]
#     p_e       : 1*barsa    -> 0.5*barsa
#     p_cap     : 5*barsa    -> 3*barsa
# - Add two new parameter rows:
#     row 29: poro     | 0.5
#     row 30: lowperm  | 20 md

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order below mirrors the order strings were introduced so the shared
# string table layout matches the original authoring session.
$ws.Range("A29").Value = "poro"
$ws.Range("B16").Value = "[75, 75]"
$ws.Range("B24").Value = "0.5*barsa"
$ws.Range("B25").Value = "3*barsa"
$ws.Range("B17").Value = "[800, 450]"
$ws.Range("A30").Value = "lowperm"
$ws.Range("B30").Value = "20 md"

$ws.Range("B29").Value = 0.5

# Update view state to match the author's final selection / scroll position
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("B29").Select()
